$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("C2").Value = 2831
$ws.Range("K2").Value = 4051
$ws.Range("K3").Value = 4158
$ws.Range("K4").Value = 835
$ws.Range("K5").Value = 295
$ws.Range("K6").Value = 4652
$ws.Range("C7").Value = 13604
$ws.Range("K7").Value = 13991

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("K2").Value = 120
$ws.Range("K5").Value = 33
$ws.Range("K7").Value = 406
$ws.Range("K8").Value = 954
$ws.Range("K9").Value = 57
$ws.Range("K11").Value = 275
$ws.Range("K16").Value = 42
$ws.Range("K18").Value = 95
$ws.Range("K19").Value = 428
$ws.Range("K20").Value = 311
$ws.Range("K23").Value = 140
$ws.Range("K25").Value = 64
$ws.Range("K27").Value = 138
$ws.Range("K29").Value = 740
$ws.Range("K31").Value = 153
$ws.Range("K33").Value = 581
$ws.Range("K36").Value = 176
$ws.Range("K37").Value = 476
$ws.Range("K42").Value = 497
$ws.Range("K43").Value = 124
$ws.Range("K44").Value = 128
$ws.Range("K47").Value = 80
$ws.Range("K48").Value = 180
$ws.Range("K49").Value = 80
$ws.Range("K51").Value = 170
$ws.Range("K52").Value = 381
$ws.Range("K53").Value = 187
$ws.Range("K55").Value = 157
$ws.Range("C63").Value = 147
$ws.Range("K63").Value = 43
$ws.Range("K65").Value = 322
$ws.Range("K68").Value = 35
$ws.Range("K71").Value = 43
$ws.Range("K73").Value = 128
$ws.Range("K74").Value = 15
$ws.Range("K75").Value = 46
$ws.Range("K76").Value = 199
$ws.Range("K77").Value = 98
$ws.Range("K78").Value = 166
$ws.Range("K79").Value = 364
$ws.Range("K83").Value = 299
$ws.Range("K84").Value = 101
$ws.Range("K85").Value = 630
$ws.Range("K88").Value = 161
$ws.Range("K89").Value = 197
$ws.Range("K90").Value = 129
$ws.Range("K91").Value = 153
$ws.Range("K92").Value = 49
$ws.Range("K93").Value = 49
$ws.Range("K94").Value = 176
$ws.Range("K95").Value = 237
$ws.Range("K99").Value = 241
$ws.Range("C101").Value = 13604
$ws.Range("K101").Value = 13991

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item(5)
$ws.Range("K2").Value = 145
$ws.Range("K3").Value = 133
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 406

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 275

# Sheet 7: Uptown
$ws = $wb.Worksheets.Item(7)
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 197

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item(8)
$ws.Range("K3").Value = 210
$ws.Range("K7").Value = 630

# Sheet 9: Little Village
$ws = $wb.Worksheets.Item(9)
$ws.Range("K3").Value = 100
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 381

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item(11)
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 187

# Sheet 12: Austin
$ws = $wb.Worksheets.Item(12)
$ws.Range("K2").Value = 268
$ws.Range("K3").Value = 284
$ws.Range("K5").Value = 27
$ws.Range("K6").Value = 321
$ws.Range("K7").Value = 954

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item(13)
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 299

# Sheet 14: Garfield Park
$ws = $wb.Worksheets.Item(14)
$ws.Range("K3").Value = 219
$ws.Range("K7").Value = 581

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item(15)
$ws.Range("K3").Value = 84
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 237

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item(16)
$ws.Range("K3").Value = 159
$ws.Range("K6").Value = 142
$ws.Range("K7").Value = 476

# Sheet 17: New City
$ws = $wb.Worksheets.Item(17)
$ws.Range("K2").Value = 95
$ws.Range("K7").Value = 322

# Sheet 18: Woodlawn
$ws = $wb.Worksheets.Item(18)
$ws.Range("K3").Value = 98
$ws.Range("K7").Value = 241

# Sheet 20: Gage Park
$ws = $wb.Worksheets.Item(20)
$ws.Range("K2").Value = 54
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 153

# Sheet 22: South Deering
$ws = $wb.Worksheets.Item(22)
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 101

# Sheet 23: Lincoln Park
$ws = $wb.Worksheets.Item(23)
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 80

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Range("K2").Value = 210
$ws.Range("K3").Value = 264
$ws.Range("K5").Value = 22
$ws.Range("K7").Value = 740

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item(26)
$ws.Range("K3").Value = 42
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 180

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item(27)
$ws.Range("K3").Value = 130
$ws.Range("K6").Value = 129
$ws.Range("K7").Value = 428

# Sheet 28: Irving Park
$ws = $wb.Worksheets.Item(28)
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 128

# Sheet 29: River North
$ws = $wb.Worksheets.Item(29)
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 199

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Range("K3").Value = 162
$ws.Range("K6").Value = 178
$ws.Range("K7").Value = 497

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item(35)
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 166

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item(36)
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 157

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item(39)
$ws.Range("K3").Value = 51
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 140

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item(40)
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 153

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item(42)
$ws.Range("K2").Value = 123
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 364

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item(44)
$ws.Range("K2").Value = 108
$ws.Range("K3").Value = 95
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 311

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item(45)
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 95

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item(47)
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 176

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item(48)
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 49

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 32
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 176

# Sheet 52: East Side
$ws = $wb.Worksheets.Item(52)
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 64

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item(53)
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 80

# Sheet 61: Avalon Park
$ws = $wb.Worksheets.Item(61)
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 57

# Sheet 62: Portage Park
$ws = $wb.Worksheets.Item(62)
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 128

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item(64)
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 120

# Sheet 66: West Elsdon
$ws = $wb.Worksheets.Item(66)
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 49

# Sheet 68: United Center
$ws = $wb.Worksheets.Item(68)
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 161

# Sheet 70: Armour Square
$ws = $wb.Worksheets.Item(70)
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 33

# Sheet 71: Edgewater
$ws = $wb.Worksheets.Item(71)
$ws.Range("K2").Value = 35
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 138

# Sheet 73: Pullman
$ws = $wb.Worksheets.Item(73)
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 46

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item(74)
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 129

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item(75)
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 170

# Sheet 76: North Park
$ws = $wb.Worksheets.Item(76)
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 35

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item(79)
$ws.Range("K2").Value = 25
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 124

# Sheet 81: Oakland
$ws = $wb.Worksheets.Item(81)
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 43

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item(84)
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 98

# Sheet 94: Bucktown
$ws = $wb.Worksheets.Item(94)
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 42

# Sheet 95: Printers Row
$ws = $wb.Worksheets.Item(95)
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 15
